# Updates the cryptos list worksheet (Sheet1) with refreshed price/volume data.
# A leading "'" forces Excel to store the value as text (matching the original
# inline-string cells) instead of auto-converting numeric-looking text to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.575.94"
$ws.Range("E2").Value = "'  -0.42%  "
$ws.Range("D3").Value = "'1.691.38"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("D5").Value = "'314.27"
$ws.Range("E5").Value = "'  -0.73%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  -0.04%  "
$ws.Range("D7").Value = "'0.3894"
$ws.Range("E7").Value = "'  -1.33%  "
$ws.Range("D8").Value = "'0.4034"
$ws.Range("E8").Value = "'  -0.61%  "
$ws.Range("D9").Value = "'1.498"
$ws.Range("E9").Value = "'  +0.68%  "
$ws.Range("D10").Value = "'1.004"
$ws.Range("E10").Value = "'  +0.07%  "
$ws.Range("D11").Value = "'52.96"
$ws.Range("E11").Value = "'  +1.58%  "
$ws.Range("D12").Value = "'0.08748"
$ws.Range("E12").Value = "'  -1.44%  "
$ws.Range("D13").Value = "'25.40"
$ws.Range("E13").Value = "'  +7.53%  "
$ws.Range("D14").Value = "'7.531"
$ws.Range("E14").Value = "'  +3.84%  "
$ws.Range("B15").Value = "'ShibaInu"
$ws.Range("C15").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001353"
$ws.Range("E15").Value = "'  +2.64%  "
$ws.Range("B16").Value = "'Chainlink"
$ws.Range("C16").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.955"
$ws.Range("E16").Value = "'  -1.15%  "
$ws.Range("D17").Value = "'1.688.20"
$ws.Range("E17").Value = "'  -0.36%  "
$ws.Range("D18").Value = "'98.60"
$ws.Range("E18").Value = "'  -1.15%  "
$ws.Range("E19").Value = "'  +1.18%  "
$ws.Range("D20").Value = "'19.93"
$ws.Range("E20").Value = "'  +1.63%  "
$ws.Range("D21").Value = "'7.297"
$ws.Range("E21").Value = "'  +4.18%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  -0.40%  "
$ws.Range("E23").Value = "'  -0.90%  "
$ws.Range("D24").Value = "'24.561.28"
$ws.Range("E24").Value = "'  -0.38%  "
$ws.Range("D25").Value = "'3.009"
$ws.Range("E25").Value = "'  -6.44%  "
$ws.Range("D26").Value = "'2.352"
$ws.Range("E26").Value = "'  -0.30%  "
$ws.Range("D27").Value = "'22.84"
$ws.Range("E27").Value = "'  +0.42%  "
$ws.Range("D28").Value = "'162.14"
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("D29").Value = "'8.755"
$ws.Range("E29").Value = "'  +14.84%  "
$ws.Range("D30").Value = "'137.15"
$ws.Range("E30").Value = "'  +0.97%  "
$ws.Range("D31").Value = "'5.218"
$ws.Range("E31").Value = "'  +0.74%  "
$ws.Range("D32").Value = "'1.872.43"
$ws.Range("E32").Value = "'  -0.42%  "
$ws.Range("E33").Value = "'  +2.58%  "
$ws.Range("D34").Value = "'7.432"
$ws.Range("E34").Value = "'  +4.53%  "
$ws.Range("D35").Value = "'1.038"
$ws.Range("E35").Value = "'  -1.98%  "
$ws.Range("D36").Value = "'1.973"
$ws.Range("E36").Value = "'  +4.49%  "
$ws.Range("D37").Value = "'0.02925"
$ws.Range("E37").Value = "'  +7.46%  "
$ws.Range("D38").Value = "'0.2749"
$ws.Range("E38").Value = "'  +0.58%  "
$ws.Range("D39").Value = "'10.78"
$ws.Range("E39").Value = "'  -4.83%  "
$ws.Range("D40").Value = "'14.30"
$ws.Range("E40").Value = "'  -1.40%  "
$ws.Range("D41").Value = "'0.09141"
$ws.Range("E41").Value = "'  -0.80%  "
$ws.Range("D42").Value = "'0.7899"
$ws.Range("E42").Value = "'  +2.99%  "
$ws.Range("D43").Value = "'1.461"
$ws.Range("E43").Value = "'  -0.57%  "
$ws.Range("E44").Value = "'  +4.93%  "
$ws.Range("D45").Value = "'0.7220"
$ws.Range("E45").Value = "'  +0.89%  "
$ws.Range("D46").Value = "'2.577"
$ws.Range("E46").Value = "'  -0.81%  "
$ws.Range("D47").Value = "'4.206"
$ws.Range("E47").Value = "'  -0.35%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "'  -0.04%  "
$ws.Range("D49").Value = "'1.336"
$ws.Range("E49").Value = "'  +1.14%  "
$ws.Range("D50").Value = "'137.90"
$ws.Range("E50").Value = "'  -1.62%  "
$ws.Range("D51").Value = "'91.15"
$ws.Range("E51").Value = "'  +0.63%  "
